$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.366.84'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '1.825.70'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.89'
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9984'
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5336'
$ws.Range("E7").Value = '  -1.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3992'
$ws.Range("E8").Value = '  +5.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07542'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.85'
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.103'
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.305'
$ws.Range("E12").Value = '  +2.34%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.626'
$ws.Range("E13").Value = '  +3.91%  '
$ws.Range("B14").Value = 'BinanceUSD'
$ws.Range("C14").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.9974'
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").Value = '1.821.41'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.61'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06580'
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.46'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9970'
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.023'
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("D23").Value = '28.390.82'
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.19'
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.075'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.80'
$ws.Range("E26").Value = '  -2.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.50'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = '2.024.24'
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.389'
$ws.Range("E29").Value = '  +1.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.29'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1098'
$ws.Range("E31").Value = '  +3.77%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.109'
$ws.Range("E32").Value = '  -2.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.681'
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.590'
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07312'
$ws.Range("E35").Value = '  +11.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2235'
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.221'
$ws.Range("E37").Value = '  +3.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02312'
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.691'
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.29'
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6221'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.193'
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.408'
$ws.Range("E43").Value = '  -2.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.41'
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.701'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5790'
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.09'
$ws.Range("E47").Value = '  -1.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.955'
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.188'
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06878'
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("E51").Value = '  -1.63%  '
